# Extend the "parallel" contingency table by two more columns (P, Q),
# matching the style of the existing header row (row 1) and recompute the
# per-row pattern (columns I, K, M, O flip, and new P/Q columns are added).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header): add P1=14, Q1=15, same style as the rest of row 1 ---
# (bold font, thin box border, centered/top-aligned -- matches the style
# already used by B1:O1)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$hdr = $ws.Range("P1:Q1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1
$hdr.Borders.Weight = 2

# --- Body rows 2..25: update I/K/M/O and add P/Q for every row ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P -> 2 (new)
    $ws.Cells.Item($r, 17).Value = 2   # Q -> 2 (new)
}
